# Update odds values for row 7 (match: Colorado Rapids vs Los Angeles Galaxy)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("U7").Value = 1.4
$ws.Range("V7").Value = 2.75
$ws.Range("AC7").Value = 21
$ws.Range("AG7").Value = 81
$ws.Range("AI7").Value = 15
$ws.Range("AR7").Value = 51
